$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin names, links, percentage strings, and
# price strings that Excel will not misinterpret as numbers)
$plainValues = @{
    'D2' = '61.982.77'
    'E2' = '  -1.48%  '
    'D3' = '3.418.55'
    'E3' = '  -0.51%  '
    'E4' = '  -0.01%  '
    'E5' = '  +0.29%  '
    'E6' = '  +5.11%  '
    'E7' = '  +0.08%  '
    'E8' = '  +1.52%  '
    'E9' = '  +3.17%  '
    'E10' = '  +0.79%  '
    'E11' = '  +3.86%  '
    'D12' = '4.003.54'
    'E12' = '  -0.52%  '
    'E13' = '  +0.81%  '
    'E14' = '  -1.09%  '
    'B15' = 'ShibaInu'
    'C15' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'E15' = '  +0.58%  '
    'B16' = 'WrappedEther'
    'C16' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D16' = '3.420.47'
    'E16' = '  -0.23%  '
    'D17' = '62.003.02'
    'E17' = '  -1.46%  '
    'E18' = '  +3.02%  '
    'E19' = '  +0.33%  '
    'E20' = '  -2.56%  '
    'E21' = '  -0.82%  '
    'E22' = '  +1.61%  '
    'E23' = '  +2.18%  '
    'E24' = '  +0.04%  '
    'D25' = '3.561.40'
    'E25' = '  -0.84%  '
    'E27' = '  -0.85%  '
    'E28' = '  +1.13%  '
    'E29' = '  -0.08%  '
    'E30' = '  +1.00%  '
    'E31' = '  -2.59%  '
    'E32' = '  -0.02%  '
    'E33' = '  -0.04%  '
    'E34' = '  +1.48%  '
    'E35' = '  +5.53%  '
    'E36' = '  +0.70%  '
    'E37' = '  -1.53%  '
    'E38' = '  +0.08%  '
    'E39' = '  -1.97%  '
    'D40' = '3.454.47'
    'E40' = '  -0.53%  '
    'E41' = '  +2.49%  '
    'E42' = '  +0.94%  '
    'E43' = '  -0.91%  '
    'E44' = '  +1.77%  '
    'E45' = '  -2.16%  '
    'E46' = '  -2.49%  '
    'D47' = '2.549.54'
    'E47' = '  -0.95%  '
    'B48' = 'InjectiveProtocol'
    'C48' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'E48' = '  +2.11%  '
    'B49' = 'Cosmos'
    'C49' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E49' = '  +0.30%  '
    'E50' = '  -2.89%  '
    'E51' = '  +0.03%  '
}

foreach ($ref in $plainValues.Keys) {
    $ws.Range($ref).Value = $plainValues[$ref]
}

# Price cells whose new text looks like a number to Excel (it would
# otherwise silently convert them to a numeric value and drop the
# original text formatting, e.g. "1.00" -> 1, "0.0000172" -> 1.72E-05).
# Force the cell to Text format just long enough to assign the literal
# string, then restore the original cell style so no visible
# formatting change is left behind.
$textForcedValues = @{
    'D5' = '579.14'
    'D6' = '153.94'
    'D9' = '8.02'
    'D10' = '0.125'
    'D14' = '28.63'
    'D15' = '0.0000172'
    'D18' = '6.57'
    'D19' = '14.42'
    'D20' = '8.95'
    'D21' = '382.17'
    'D22' = '0.569'
    'D23' = '76.06'
    'D27' = '0.178'
    'D31' = '7.88'
    'D32' = '1.00'
    'D33' = '23.26'
    'D37' = '6.96'
    'D41' = '0.0783'
    'D42' = '42.73'
    'D46' = '1.17'
    'D48' = '23.13'
    'D49' = '6.81'
    'D50' = '2.20'
}

foreach ($ref in $textForcedValues.Keys) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedValues[$ref]
    $cell.Style = $origStyle
}
